$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1208
$ws.Range("I33").Value = 196
$ws.Range("J33").Value = 2220
$ws.Range("K33").Value = 196
$ws.Range("L33").Value = 2220
$ws.Range("M33").Value = 33
$ws.Range("N33").Value = -2678

$ws.Range("H40").Value = 1199
$ws.Range("I40").Value = 1199
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1199
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1024
$ws.Range("N40").ClearContents()

$ws.Range("H64").Value = 8030.625

$ws.Range("H67").Value = 8030.625

$ws.Range("H70").Value = 41672070
$ws.Range("I70").Value = 3596.375
$ws.Range("J70").Value = 62506310
$ws.Range("K70").Value = 10789.125
$ws.Range("L70").Value = 187518930
$ws.Range("M70").Value = -10519.125
$ws.Range("N70").Value = -187519470

$ws.Range("H73").Value = 41672070
$ws.Range("I73").Value = 3596.375
$ws.Range("J73").Value = 62506310
$ws.Range("K73").Value = 10789.125
$ws.Range("L73").Value = 187518930
$ws.Range("M73").Value = -9853.125
$ws.Range("N73").Value = -187520802

$ws.Range("H129").Value = 1543.4445
$ws.Range("I129").Value = 852.3333
$ws.Range("K129").Value = 2556.9999
$ws.Range("M129").Value = 2443.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5078.0713
$ws.Range("I32").Value = 5286.9165
$ws.Range("J32").Value = 3825
$ws.Range("K32").Value = 5286.9165
$ws.Range("L32").Value = 3825
$ws.Range("M32").Value = -4999.9165
$ws.Range("N32").Value = -4399

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 1684.6
$ws.Range("I33").Value = 1349.75
$ws.Range("J33").Value = 3024
$ws.Range("K33").Value = 1349.75
$ws.Range("L33").Value = 3024
$ws.Range("M33").Value = -1013.75
$ws.Range("N33").Value = -3696

$ws.Range("H36").Value = 977.5
$ws.Range("I36").Value = 977.5
$ws.Range("K36").Value = 977.5
$ws.Range("M36").Value = -443.5

$ws.Range("H134").Value = 2364.074
$ws.Range("I134").Value = 2344.5264
$ws.Range("K134").Value = 7033.5792
$ws.Range("M134").Value = -4498.5792

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 726114.7
$ws.Range("I6").Value = 1004260.2
$ws.Range("K6").Value = 1004260.2
$ws.Range("M6").Value = -1004147.2

$ws.Range("H62").Value = 5190.7334
$ws.Range("I62").Value = 6240.5
$ws.Range("J62").Value = 3091.2
$ws.Range("K62").Value = 6240.5
$ws.Range("L62").Value = 3091.2
$ws.Range("M62").Value = -5616.5
$ws.Range("N62").Value = -4339.2

$ws.Range("H65").Value = 5190.7334
$ws.Range("I65").Value = 6240.5
$ws.Range("J65").Value = 3091.2
$ws.Range("K65").Value = 31202.5
$ws.Range("L65").Value = 15456
$ws.Range("M65").Value = -28082.5
$ws.Range("N65").Value = -21696

$ws.Range("H122").Value = 1706.3846
$ws.Range("I122").Value = 1364.3334
$ws.Range("J122").Value = 1999.5714
$ws.Range("K122").Value = 4093.0002
$ws.Range("L122").Value = 5998.7142
$ws.Range("M122").Value = -1643.0002
$ws.Range("N122").Value = -10898.7142

$ws.Range("H134").Value = 3750.1333
$ws.Range("I134").Value = 3408.75
$ws.Range("K134").Value = 10226.25
$ws.Range("M134").Value = -7691.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 679.9167
$ws.Range("I5").Value = 679.9167
$ws.Range("K5").Value = 2039.7501
$ws.Range("M5").Value = -1927.7501

$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 1000
$ws.Range("K8").Value = 3000
$ws.Range("M8").Value = -2861

$ws.Range("H25").Value = 253.5
$ws.Range("J25").Value = 174.25
$ws.Range("L25").Value = 522.75
$ws.Range("N25").Value = -860.75

$ws.Range("H30").Value = 253.5
$ws.Range("J30").Value = 174.25
$ws.Range("L30").Value = 522.75
$ws.Range("N30").Value = -726.75

$ws.Range("H32").Value = 1150
$ws.Range("J32").Value = 1150
$ws.Range("L32").Value = 3450
$ws.Range("N32").Value = -4016

$ws.Range("H33").Value = 319.91666
$ws.Range("J33").Value = 376.55554
$ws.Range("L33").Value = 2259.33324
$ws.Range("N33").Value = -2825.33324

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H36").Value = 1250
$ws.Range("I36").Value = 1250
$ws.Range("K36").Value = 3750
$ws.Range("M36").Value = -3581

$ws.Range("H38").Value = 876.8333
$ws.Range("I38").Value = 27
$ws.Range("J38").Value = 1301.75
$ws.Range("K38").Value = 81
$ws.Range("L38").Value = 3905.25
$ws.Range("M38").Value = 266
$ws.Range("N38").Value = -4599.25

$ws.Range("H60").Value = 2105
$ws.Range("I60").Value = 2105
$ws.Range("K60").Value = 6315
$ws.Range("M60").Value = -6064

$ws.Range("H68").Value = 6599.6
$ws.Range("J68").Value = 6599.6
$ws.Range("L68").Value = 19798.8
$ws.Range("N68").Value = -21420.8

$ws.Range("H71").Value = 6599.6
$ws.Range("J71").Value = 6599.6
$ws.Range("L71").Value = 59396.4
$ws.Range("N71").Value = -67508.39999999999

$ws.Range("H135").Value = 679.9167
$ws.Range("I135").Value = 679.9167
$ws.Range("K135").Value = 6119.2503
$ws.Range("M135").Value = -3584.2503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4387.3477
$ws.Range("I2").Value = 7736.077
$ws.Range("J2").Value = 34
$ws.Range("K2").Value = 7736.077
$ws.Range("L2").Value = 34
$ws.Range("M2").Value = -7623.077
$ws.Range("N2").Value = -260

$ws.Range("H59").Value = 14646.363
$ws.Range("I59").Value = 8000
$ws.Range("J59").Value = 44555
$ws.Range("K59").Value = 8000
$ws.Range("L59").Value = 44555
$ws.Range("M59").Value = -7417
$ws.Range("N59").Value = -45721

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 881.44446
$ws.Range("I16").Value = 866.625
$ws.Range("K16").Value = 866.625
$ws.Range("M16").Value = -696.625

$ws.Range("H22").Value = 3619.5881
$ws.Range("I22").Value = 3403.9
$ws.Range("J22").Value = 3927.7144
$ws.Range("K22").Value = 3403.9
$ws.Range("L22").Value = 3927.7144
$ws.Range("M22").Value = -3108.9
$ws.Range("N22").Value = -4517.7144

$ws.Range("H27").Value = 3619.5881
$ws.Range("I27").Value = 3403.9
$ws.Range("J27").Value = 3927.7144
$ws.Range("K27").Value = 3403.9
$ws.Range("L27").Value = 3927.7144
$ws.Range("M27").Value = -3296.9
$ws.Range("N27").Value = -4141.7144

$ws.Range("H68").Value = 6657.1055
$ws.Range("J68").Value = 7875
$ws.Range("L68").Value = 7875
$ws.Range("N68").Value = -9373

$ws.Range("H71").Value = 6657.1055
$ws.Range("J71").Value = 7875
$ws.Range("L71").Value = 39375
$ws.Range("N71").Value = -46863

$ws.Range("H74").Value = 57000
$ws.Range("I74").Value = 54000
$ws.Range("K74").Value = 54000
$ws.Range("M74").Value = -53002

$ws.Range("H77").Value = 57000
$ws.Range("I77").Value = 54000
$ws.Range("K77").Value = 162000
$ws.Range("M77").Value = -157008

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1313.2667
$ws.Range("I14").Value = 1180
$ws.Range("K14").Value = 1180
$ws.Range("M14").Value = -1012

$ws.Range("H37").Value = 25000
$ws.Range("J37").Value = 25000
$ws.Range("L37").Value = 25000
$ws.Range("N37").Value = -25406

$ws.Range("H62").Value = 10512.462
$ws.Range("J62").Value = 10692.96
$ws.Range("L62").Value = 10692.96
$ws.Range("N62").Value = -11940.96

$ws.Range("H65").Value = 10512.462
$ws.Range("J65").Value = 10692.96
$ws.Range("L65").Value = 53464.8
$ws.Range("N65").Value = -59704.8

$ws.Range("H81").Value = 2206.4614
$ws.Range("I81").Value = 1607.6364
$ws.Range("J81").Value = 5500
$ws.Range("K81").Value = 3215.2728
$ws.Range("L81").Value = 11000
$ws.Range("M81").Value = -2154.2728
$ws.Range("N81").Value = -13122

$ws.Range("H84").Value = 2206.4614
$ws.Range("I84").Value = 1607.6364
$ws.Range("J84").Value = 5500
$ws.Range("K84").Value = 16076.364
$ws.Range("L84").Value = 55000
$ws.Range("M84").Value = -10772.364
$ws.Range("N84").Value = -65608

$ws.Range("H96").Value = 3250
$ws.Range("I96").Value = 3000
$ws.Range("J96").Value = 3500
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 3500
$ws.Range("M96").Value = -1627
$ws.Range("N96").Value = -6246
